$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "lambda" column header (Argon block)
$ws.Range("J4").Value = "lambda"

# New lambda readings for the Argon data rows
$ws.Range("J5").Value = 0.270406
$ws.Range("J6").Value = 0.240883

# "avg" label plus average / standard-error-of-the-mean formulas
$ws.Range("I7").Value = "avg"
$ws.Range("J7").Formula = "=AVERAGE(J5:J6)"
$ws.Range("K7").Formula = "=STDEV(J5:J6)/SQRT(2)"

# Touch every populated row's height so it gets persisted as an explicit
# (custom) height, matching the saved workbook's row metadata.
$rows = @(1,3,4,5,6,7,8,9,10,12,13,14,15,16,17,18,19,21,22,23,24,25,26,27,28)
foreach ($r in $rows) {
    $row = $ws.Rows.Item($r)
    $row.RowHeight = $row.RowHeight
}

# Match the author's final selection
$null = $ws.Range("K8").Select()
